$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.232.58"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "2.784.85"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'589.14"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Value = "'161.34"
$ws.Range("E6").Value = "  +7.82%  "
$ws.Range("D7").Value = "'0.620"
$ws.Range("E7").Value = "  +2.08%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("E11").Value = "  +3.27%  "
$ws.Range("D12").Value = "'0.160"
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("D13").Value = "3.279.93"
$ws.Range("E13").Value = "  +2.76%  "
$ws.Range("D14").Value = "'27.67"
$ws.Range("E14").Value = "  +4.26%  "
$ws.Range("D15").Value = "64.145.04"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("E16").Value = "  +6.45%  "
$ws.Range("D17").Value = "2.787.62"
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("D18").Value = "'12.43"
$ws.Range("E18").Value = "  +4.22%  "
$ws.Range("D19").Value = "'5.09"
$ws.Range("E19").Value = "  +4.49%  "
$ws.Range("D20").Value = "'367.87"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").Value = "'7.09"
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("E22").Value = "  +8.27%  "
$ws.Range("E23").Value = "  +0.78%  "
$ws.Range("D24").Value = "'67.66"
$ws.Range("E24").Value = "  +3.40%  "
$ws.Range("E25").Value = "  +6.91%  "
$ws.Range("D26").Value = "'8.90"
$ws.Range("E26").Value = "  +4.16%  "
$ws.Range("D27").Value = "0.0₃0974"
$ws.Range("E27").Value = "  +13.71%  "
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("E29").Value = "  +1.35%  "
$ws.Range("E30").Value = "  +3.22%  "
$ws.Range("D31").Value = "'1.27"
$ws.Range("E31").Value = "  +6.59%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'5.23"
$ws.Range("E32").Value = "  +9.81%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "'172.08"
$ws.Range("E33").Value = "  +1.38%  "
$ws.Range("D34").Value = "'20.92"
$ws.Range("E34").Value = "  +1.96%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("E36").Value = "  +5.74%  "
$ws.Range("E37").Value = "  +2.29%  "
$ws.Range("E38").Value = "  +2.62%  "
$ws.Range("E39").Value = "  +0.87%  "
$ws.Range("E40").Value = "  +12.43%  "
$ws.Range("D41").Value = "'341.61"
$ws.Range("E41").Value = "  -2.70%  "
$ws.Range("D42").Value = "'40.32"
$ws.Range("E42").Value = "  +2.65%  "
$ws.Range("D43").Value = "'22.58"
$ws.Range("E43").Value = "  +4.65%  "
$ws.Range("D44").Value = "'22.61"
$ws.Range("E44").Value = "  +4.31%  "
$ws.Range("E45").Value = "  +3.63%  "
$ws.Range("D46").Value = "'0.654"
$ws.Range("E46").Value = "  +2.24%  "
$ws.Range("D47").Value = "'0.0264"
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("D48").Value = "'138.83"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("D50").Value = "2.178.98"
$ws.Range("E50").Value = "  +1.69%  "
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  +0.40%  "
